# Correção nos dados e inicio da analise PNAD 2009
#
# The sheet had two "section header" rows that only carried a label in
# column A ("situação do domicílio" and "grandes regiões e unidades da
# federação") with no data beside them. Those are removed entirely
# (shifting everything below them up), and the now-redundant "unnamed:
# 1_level_1" / "unnamed: 5_level_1" sub-header labels in row 2 are
# corrected to read "total" (matching column C / column F already
# being "total").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the row-2 sub-headers that pandas had left as "unnamed: x_level_1"
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"

# Remove the "situação do domicílio" section-header row (row 5).
$ws.Rows.Item(5).Delete()

# After the deletion above, the "grandes regiões e unidades da
# federação" section-header row shifted up from row 8 to row 7.
$ws.Rows.Item(7).Delete()
